$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) number format, used to restore
# style on cells whose text-like numeric values Excel would otherwise
# auto-convert to numbers.
$defaultStyle = $ws.Range("D4").Style

# Row 2
$ws.Range("D2").Value = "64.221.67"
$ws.Range("E2").Value = "  -0.42%  "

# Row 3
$ws.Range("D3").Value = "3.139.09"
$ws.Range("E3").Value = "  -1.26%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.90"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.22%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.99"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -3.67%  "

# Row 7
$ws.Range("E7").Value = "  -0.25%  "

# Row 8
$ws.Range("E8").Value = "  -5.78%  "

# Row 9
$ws.Range("D9").Value = "3.151.51"
$ws.Range("E9").Value = "  -0.89%  "

# Row 10
$ws.Range("E10").Value = "  -3.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.58"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -3.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  -0.22%  "

# Row 13
$ws.Range("D13").Value = "3.686.18"

# Row 14
$ws.Range("E14").Value = "  -0.55%  "

# Row 15
$ws.Range("D15").Value = "64.273.07"
$ws.Range("E15").Value = "  -0.43%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.09"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  -1.18%  "

# Row 17
$ws.Range("D17").Value = "3.146.82"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18
$ws.Range("E18").Value = "  -3.63%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "400.54"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  -4.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.23"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -2.51%  "

# Row 21
$ws.Range("E21").Value = "  -3.07%  "

# Row 22
$ws.Range("E22").Value = "  -0.44%  "

# Row 23
$ws.Range("E23").Value = "  +0.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.08"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -2.90%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.483"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -1.11%  "

# Row 26
$ws.Range("E26").Value = "  -4.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000101"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -5.04%  "

# Row 28
$ws.Range("E28").Value = "  -1.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -0.38%  "

# Row 30
$ws.Range("E30").Value = "  +0.05%  "

# Row 31
$ws.Range("E31").Value = "  -1.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.17"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  -2.51%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.44"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +1.76%  "

# Row 34
$ws.Range("E34").Value = "  -1.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.81"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -5.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -2.91%  "

# Row 37
$ws.Range("E37").Value = "  -2.57%  "

# Row 38
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.68"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -2.15%  "

# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.648.98"
$ws.Range("E39").Value = "  -2.82%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.59"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -2.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.06"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -2.90%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.30"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  -2.25%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0611"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -2.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.47"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -2.75%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0254"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -3.98%  "

# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "285.01"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -3.13%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.00"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -3.00%  "

# Row 49
$ws.Range("E49").Value = "  -0.32%  "

# Row 50
$ws.Range("E50").Value = "  -1.52%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.46"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -0.05%  "
